$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.746.42"
$ws.Range('E2').Value = '  -3.60%  '
$ws.Range('D3').Value = "'2.903.59"
$ws.Range('E3').Value = '  -4.25%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'585.43"
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').Value = "'144.82"
$ws.Range('E6').Value = '  -5.92%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = "'0.503"
$ws.Range('E8').Value = '  -2.64%  '
$ws.Range('D9').Value = "'2.902.35"
$ws.Range('E9').Value = '  -4.18%  '
$ws.Range('D10').Value = "'6.66"
$ws.Range('E10').Value = '  -2.72%  '
$ws.Range('E11').Value = '  -4.67%  '
$ws.Range('E12').Value = '  -3.81%  '
$ws.Range('D14').Value = "'33.42"
$ws.Range('E14').Value = '  -6.37%  '
$ws.Range('D16').Value = "'3.384.20"
$ws.Range('E16').Value = '  -4.24%  '
$ws.Range('D17').Value = "'60.729.86"
$ws.Range('D18').Value = "'6.70"
$ws.Range('E18').Value = '  -5.37%  '
$ws.Range('D19').Value = "'2.902.82"
$ws.Range('E19').Value = '  -4.27%  '
$ws.Range('D20').Value = "'428.26"
$ws.Range('E20').Value = '  -5.36%  '
$ws.Range('D21').Value = "'13.57"
$ws.Range('E21').Value = '  -4.84%  '
$ws.Range('E22').Value = '  -2.43%  '
$ws.Range('D23').Value = "'7.11"
$ws.Range('E23').Value = '  -5.28%  '
$ws.Range('D24').Value = "'80.64"
$ws.Range('E24').Value = '  -2.96%  '
$ws.Range('E25').Value = '  -3.28%  '
$ws.Range('E26').Value = '  -2.43%  '
$ws.Range('D27').Value = "'11.92"
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').Value = "'7.28"
$ws.Range('E29').Value = '  -2.73%  '
$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').Value = "'1.00"
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('E32').Value = '  -3.34%  '
$ws.Range('D33').Value = "'26.45"
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('D35').Value = "'0.0₃0874"
$ws.Range('E35').Value = '  +1.71%  '
$ws.Range('E36').Value = '  -3.18%  '
$ws.Range('E37').Value = '  -5.28%  '
$ws.Range('E38').Value = '  -4.13%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = "'49.52"
$ws.Range('E39').Value = '  -1.97%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = "'0.125"
$ws.Range('E40').Value = '  -3.93%  '
$ws.Range('E41').Value = '  -4.61%  '
$ws.Range('D42').Value = "'8.61"
$ws.Range('E42').Value = '  -5.66%  '
$ws.Range('D43').Value = "'0.297"
$ws.Range('E43').Value = '  -3.51%  '
$ws.Range('E44').Value = '  -5.95%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = "'0.0351"
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = "'378.19"
$ws.Range('E46').Value = '  -2.89%  '
$ws.Range('D47').Value = "'2.698.17"
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('D48').Value = "'132.65"
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('D50').Value = "'24.25"
$ws.Range('E50').Value = '  -3.12%  '
$ws.Range('E51').Value = '  -2.52%  '
